$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.714.26"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "3.579.56"

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "604.15"
$ws.Range("E5").Value = "  +1.24%  "

$ws.Range("D6").Value = "137.04"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").Value = "3.579.24"
$ws.Range("E7").Value = "  +1.39%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  +0.55%  "

$ws.Range("D11").Value = "7.21"
$ws.Range("E11").Value = "  +5.96%  "

$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("D13").Value = "4.193.12"
$ws.Range("E13").Value = "  +1.38%  "

$ws.Range("D14").Value = "28.16"
$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").Value = "3.583.71"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "65.771.44"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "10.07"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("D20").Value = "14.65"
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("D21").Value = "5.87"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").Value = "394.29"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("D24").Value = "3.723.73"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").Value = "74.10"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").Value = "1.00"

$ws.Range("E27").Value = "  +2.05%  "

$ws.Range("E28").Value = "  +4.45%  "

$ws.Range("E29").Value = "  +26.04%  "

$ws.Range("D30").Value = "2.35"
$ws.Range("E30").Value = "  +2.63%  "

$ws.Range("D31").Value = "8.58"
$ws.Range("E31").Value = "  +5.45%  "

$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("D33").Value = "3.584.06"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").Value = "24.48"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("D35").Value = "0.148"
$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("D37").Value = "5.35"
$ws.Range("E37").Value = "  +7.12%  "

$ws.Range("E38").Value = "  +5.22%  "

$ws.Range("D39").Value = "7.04"
$ws.Range("E39").Value = "  +0.98%  "

$ws.Range("D40").Value = "167.42"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("E41").Value = "  +4.02%  "

$ws.Range("D42").Value = "0.838"
$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").Value = "26.93"
$ws.Range("E43").Value = "  +3.46%  "

$ws.Range("D44").Value = "1.27"
$ws.Range("E44").Value = "  +7.53%  "

$ws.Range("E46").Value = "  +2.41%  "

$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").Value = "  +1.73%  "

$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").Value = "2.455.98"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0270"
$ws.Range("E51").Value = "  +4.50%  "

